$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '71.844.06'
$ws.Range("E2").Value = '  +3.70%  '
$ws.Range("D3").Value = '3.693.69'
$ws.Range("E3").Value = '  +8.34%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("D7").Value = '3.683.34'
$ws.Range("E7").Value = '  +8.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.623'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.18%  '
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("E11").Value = '  +4.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '50.21'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.51%  '
$ws.Range("E13").Value = '  +2.24%  '
$ws.Range("D14").Value = '4.290.23'
$ws.Range("E14").Value = '  +8.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '687.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '9.05'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.67%  '
$ws.Range("D17").Value = '3.692.89'
$ws.Range("E17").Value = '  +8.45%  '
$ws.Range("D18").Value = '71.909.00'
$ws.Range("E18").Value = '  +3.63%  '
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.47%  '
$ws.Range("E22").Value = '  +3.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +17.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '104.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.54%  '
$ws.Range("E26").Value = '  +3.97%  '
$ws.Range("E27").Value = '  +5.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.26'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.82%  '
$ws.Range("E31").Value = '  +6.83%  '
$ws.Range("E32").Value = '  +13.33%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.89%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '569.73'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.16%  '
$ws.Range("E35").Value = '  +4.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.56'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.75%  '
$ws.Range("D37").Value = '3.809.79'
$ws.Range("E37").Value = '  +5.42%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D40").Value = '0.0₃0782'
$ws.Range("E40").Value = '  +5.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("E42").Value = '  +6.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0467'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.51%  '
$ws.Range("E45").Value = '  +5.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.92'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.34%  '
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("E48").Value = '  +4.27%  '
$ws.Range("E49").Value = '  +3.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '134.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.40%  '
